# Remove the "Role" column (E) contents from the staff sheet.
# Only "Name", "Email", "Faculty" and "Password" are kept as persistent
# data saved outside the app; the Role/committeeOf values were never
# meant to be stored here, so clear them out (column stays, just empty).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("staff")

# Select the whole column E (mirrors clicking the column header) and clear
# its contents - this keeps the column itself (and E1's formatting) in
# place instead of shifting D/C/etc left like a full column delete would.
$ws.Range("E1:E1048576").Select()
$ws.Range("E1:E1048576").ClearContents()
